# Automatische test-sync: 2025-08-03 23:33:50
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log row (row 53)
$ws.Cells.Item(53, 1).Value = "Is er al nieuws?"
$ws.Cells.Item(53, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(53, 3).Value = "Testmail #4: Is er al nieuws?"
$ws.Cells.Item(53, 4).Value = "Overig"
$ws.Cells.Item(53, 5).Value = "Bedankt, we hebben dit doorgestuurd naar support@bedrijf.nl."
$ws.Cells.Item(53, 6).Value = "2025-08-03 23:33:43"
$ws.Cells.Item(53, 7).Value = "Ja"
$ws.Cells.Item(53, 8).Value = "Ja"
$ws.Cells.Item(53, 9).Value = "Nee"
$ws.Cells.Item(53, 10).Value = "Nee"

# Extend the conditional-formatting ranges from row 52 to row 53
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $ws.Range($col + "2:" + $col + "52")
    $newRange = $ws.Range($col + "2:" + $col + "53")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the Dashboard "Overig" tally to reflect the new row
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Cells.Item(3, 2).Value = 13
